$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins / Losses / Ties) in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, centered, bordered) used by the
# existing header row by copying the format from A1 onto the new cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins=79, Losses=83, Ties=0) for every
# player row (2 through 44).
$ws.Range("AD2:AD44").Value = 79
$ws.Range("AE2:AE44").Value = 83
$ws.Range("AF2:AF44").Value = 0
